$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New / updated tutorial text blocks
$s20 = 'Congratulations, you made it through the first section of the map and proved that you are in fact a carbon based life form.  Let’s keep expanding our primitive minds with the next command.  The {move forward until unable} button is a button that can be very useful to cover long strait distances without having to push {move forward} many times over.  The {move forward until unable} button will have Rufus continue tCongratulations, you made it through the first section of the map and proved that you are in fact a carbon based life form.  Let’s keep expanding our primitive minds with the next command.  The {move forward until unable} button is a button that can be very useful to cover long strait distances without having to push {move forward} many times over.  The {move forward until unable} button will have Rufus continue to travel in the direction he is facing at the time of the buttons execution until he dies or is no longer able to move forward because something is obstructing his path; such as a wall or door or edge of the map (yes the world is flat but no you will not sail off the edge of the world).  This tool is useful but beware of over use as its byte cost is high.  If you find yourself out of bytes you may want to go back and take a few of these out of your command list.  For now though, go nuts and enjoy yourself. '
$s21 = 'Congratulations!!!  While that was no Everest don’t let anyone take that moment away from.  You too are  a unique and special snowflake.  In front of you there are three new obstacles; a gap, a duct and an electric tile.  To get past these obstacles you will need to use the {jump} button and the {crouch} button.  Both gaps and electric tiles can be jumped over by using the {jump} button.  You can even jump from one raised tile to another raised tile of the same height.  The duct, in between the two jump obstacles, must be ducked under for Rufus to continue.  For this you will use the {crouch} button.  So jump, limbo, and jump your way through to the next section.  {End Movie 05}'
$s22 = 'The next three tiles in front of you can all need to be activated by Rufus.  You can do this by use the {activate} button each time Rufus comes in contact with them. The first obstacle we see is a door that needs to be opened before Rufus can walk through it.  Doors can be on any side of a tile so it is important that you are facing the same direction as the door in order to activate and open it up.  '
$s23 = 'The next tile {reprogrammable square} is a Reprogrammable Square.  This tile operates in a similar way that these tutorial stop spots do in that it resets this to be Rufus’ new starting location should you die or need to go back and reset your current solution.  It also refunds any used bytes up to this point in your solution.  In order to access this tiles feature you will need to activate it while standing on them or facing them if they are on the edge of the tile. {top right reprogrammable square}  Reprogrammable squares can be activated as many times as you like within the level.  Making good use of them can help you increase your score by keeping your used bytes in the solution to a minimum.  '
$s24 = 'The {center switch tile} just after the Reprogrammable Square must also be activated.  This switch tile can be used for another number of different things in the game.  It has the ability to unblock paths by taking down walls or turning off electric tiles.  In addition it can take and inactive End Square {inactive end tile} and turn it back to active.  It even has the ability to make hidden portions of the map appear, like other switches, which you may need in order to get Rufus through the map.  The switch we see here is in the center of the tile which means no matter which direction Rufus is facing he can activate it.  Switches can also be on the edge of the tiles {upper left switch} and you must be facing them in order to activate them.  Take your new found power and go forth and concur.  {end movie 06}'
$s25 = 'Every time you want out….they keep pulling you back in!!!  At this point I believe you are ready to face the trials, which for this game consists of Subroutines.  As stated earlier you have a limited amount of bytes to use to navigate Rufus through each map.  You may notice parts of the map, like the one ahead, in which there is a repeating pattern to them.  In this case it is move forward three tiles and turn left or right.  Rather than put those commands in over and over again, which can be very byte expensive, we can put the commands to be repeated into the Sub1 and Sub2 section.   Why do this you ask with a puzzled look on your face?  Because I said so…okay there is more to it than that.  When you place the commands into the one of the two Sub’s it will cost you the original byte cost of each command that is placed.  But in the Main portion of the command screen when you want to execute the entire portion of code you placed into the Sub’s it will only cost you two bytes each time you use it.'
$s26 = 'For this section of code click on Sub1 on the main screen and put either {Move Forward} three times followed by {turn left} or {move forward till unable} followed by {turn left}.  Then click Sub2 and put the same thing except instead of using the {turn left} button use the {turn right} button.  Then in the Main command window use {Sub1} and {Sub2} as they match up with the map.  Give it a shot. {end movie 07}'
$s27 = 'Congratulations!!! No really this time I mean it.  Subroutines are not an easy concept to understand so give yourself a pat on the back. Just make sure you don’t choke on your gum because if you do we will lose our newfound respect for you.   Ahead of you are a group of {ice tile}; these are ice tiles.  Ice tiles are unique in that they are slippery.  Once you pick a direction on ice you only have to hit {move forward} one and Rufus will continue to slide until he hits a non-ice tile or hits an obstruction of some kind like a wall, door or edge of the map (again yes the world is flat but no you will not sail off the edge).  '
$s28 = 'The tile beyond this {teleport tile} is a teleport tile.  Hmm I wonder what this tile does.  Yes you guessed it, it teleports Rufus to another location of the map, you may now take your Fields Medal. {end movie 08}'
$s29 = 'Well you made!!! What a long strange trip it’s been.  All that’s left is to move forward on the End Square and claim this map in the name of whatever country you represent.  I hope you have enjoyed this tutorial and remember if you forget anything that has been discussed you can always click the {help button} to answer any question while you are in game.  Enjoy the game and thanks for playing!!!'

# Row 17: replace outdated 'move forward until unable' text with the corrected version
$ws.Range("C17").Value = $s20
$ws.Rows.Item(17).RowHeight = 261

# Row 16 height tweak (content unchanged)
$ws.Rows.Item(16).RowHeight = 89.25

# New rows 20-28 appended after existing row 19
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = $s21
$ws.Rows.Item(20).RowHeight = 135

$ws.Range("A21").Value = 6
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = $s22
$ws.Rows.Item(21).RowHeight = 75

$ws.Range("A22").Value = 6
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = $s23
$ws.Rows.Item(22).RowHeight = 135

$ws.Range("A23").Value = 6
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = $s24
$ws.Rows.Item(23).RowHeight = 165

$ws.Range("A24").Value = 7
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = $s25
$ws.Rows.Item(24).RowHeight = 201.75

$ws.Range("A25").Value = 7
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = $s26
$ws.Rows.Item(25).RowHeight = 90

$ws.Range("A26").Value = 8
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = $s27
$ws.Rows.Item(26).RowHeight = 120

$ws.Range("A27").Value = 8
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = $s28
$ws.Rows.Item(27).RowHeight = 45

$ws.Range("A28").Value = 9
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = $s29
$ws.Rows.Item(28).RowHeight = 81

# Update view: scroll so row 13 is at top and select D28 (matches author's final cursor position)
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("D28").Select() | Out-Null
